$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns R and S ---
$ws.Range("Q1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R1").Value = "3-jul"
$ws.Range("S1").Value = "4-jul"

# --- Data rows 2-11: new columns R and S, copy number format/style from column Q ---
$ws.Range("Q2:Q11").Copy()
$ws.Range("R2:S11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rData = @{
  2  = @(14, 13)
  3  = @(19, 21)
  4  = @(10, 10)
  5  = @(13, 12)
  6  = @(11, 14)
  7  = @(18, 18)
  8  = @(8, 10)
  9  = @(23, 23)
  10 = @(16, 23)
  11 = @(8, 3)
}

foreach ($row in $rData.Keys) {
  $vals = $rData[$row]
  $ws.Range("R$row").Value = $vals[0]
  $ws.Range("S$row").Value = $vals[1]
}

# --- Sheet view changes: drop topLeftCell, move selection ---
$ws.Range("P14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
